$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Fix the typo in the Kyrgyz header text in A1 ("Коопсуз ... жоктугунана" ->
#    "Коопсуздук ... жоктугунан"). Writing the new text directly lets the
#    engine retire the old shared-string entry and append the corrected one.
$ws.Range("A1").Value = "3.9.2 Коопсуздук суунун, коопсуздук санитариянын жана гигиенанын жоктугунан болгон өлүм"

# 2. Add the new 2022 data column (S) to the table, copying the number
#    formatting/alignment/borders from the existing 2021 column (R) and then
#    filling in the figures for every data row.
$ws.Range("R4:R14").Copy()
$ws.Range("S4:S14").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("S4").Value = 2022
$ws.Range("S5").Value = 1.2
$ws.Range("S6").Value = 2.7
$ws.Range("S7").Value = 0.9
$ws.Range("S8").Value = 0.4
$ws.Range("S9").Value = 0.7
$ws.Range("S10").Value = 0.9
$ws.Range("S11").Value = 1.1
$ws.Range("S12").Value = 2.7
$ws.Range("S13").Value = 0.4
$ws.Range("S14").Value = 0.6

# 3. Reset the selection back to A1 so the saved sheet view no longer points
#    at the old S17 selection.
$ws.Range("A1").Select()
